$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- CARBON and NITROGEN section (row 33 header) gains a "# records" count ---
$ws.Range("C33").Value = 2189
$ws.Range("C33").NumberFormat = "#,##0"
$ws.Range("D33").Value = "# records"

# --- Insert a new row for the XGBoost model result (95.8) right after the
#     existing CARBON+NITROGEN model rows (33-38), shifting everything from
#     the old row 40 onward down by one row ---
$ws.Rows.Item(39).Insert()
$ws.Range("A39").Value = "XGBoost"
$ws.Range("B39").Value = 95.8

# --- Re-apply the sort on the "C, Si, N" section (now rows 42:46) so the
#     cached sortState reference shifts along with the data ---
$sortRange = $ws.Range("A42:B46")
$key = $ws.Range("B42:B46")
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($key)
$ws.Sort.SetRange($sortRange)
$ws.Sort.Header = -4142
$ws.Sort.Apply()

# --- Update the view state to match the post-edit selection/scroll ---
$ws.Range("D33").Select()
